$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 70, shifting existing rows 70-85 down to 71-86
$ws.Rows(70).Insert()

# Fill the new row 70 with the new "Tuna" price entry
$ws.Cells.Item(70, 1).Value = 10
$ws.Cells.Item(70, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(70, 3).Value = "La Araucanía"
$ws.Cells.Item(70, 4).Value = 45015
$ws.Cells.Item(70, 5).Value = 9
$ws.Cells.Item(70, 6).Value = "Fruta"
$ws.Cells.Item(70, 7).Value = 100107
$ws.Cells.Item(70, 8).Value = "Otros"
$ws.Cells.Item(70, 9).Value = 100107011
$ws.Cells.Item(70, 10).Value = "Tuna"
$ws.Cells.Item(70, 11).Value = "Sin especificar"
$ws.Cells.Item(70, 12).Value = "Primera"
$ws.Cells.Item(70, 13).Value = 85
$ws.Cells.Item(70, 14).Value = 16000
$ws.Cells.Item(70, 15).Value = 17000
$ws.Cells.Item(70, 16).Value = 16235
$ws.Cells.Item(70, 17).Value = "$/caja 16 kilos"
$ws.Cells.Item(70, 18).Value = "Provincia de Los Andes"
$ws.Cells.Item(70, 19).Value = 1015
$ws.Cells.Item(70, 20).Value = 16
